$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4").Value = 0.001933813095092773
$ws.Range("O6").Value = 0.0006279945373535156
$ws.Range("O7").Value = 0.0005905628204345703
$ws.Range("O12").Value = 0
$ws.Range("O14").Value = 0.149749755859375
$ws.Range("O15").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("O17").Value = 0.006505727767944336
$ws.Range("O18").Value = 0
$ws.Range("O22").Value = 0.001549243927001953
$ws.Range("O23").Value = 0
$ws.Range("O24").Value = 0.001001358032226562
$ws.Range("O25").Value = 0.1840081214904785
$ws.Range("O26").Value = 0.1237947940826416
$ws.Range("O27").Value = 0.01694273948669434
$ws.Range("O29").Value = 0.008513212203979492
$ws.Range("O31").Value = 0
$ws.Range("O32").Value = 0.001675844192504883
$ws.Range("O33").Value = 0
$ws.Range("O34").Value = 7.887200117111206
$ws.Range("O35").Value = 0.01630020141601562
$ws.Range("O36").Value = 0.09285497665405273
$ws.Range("O37").Value = 0.2677047252655029
$ws.Range("O38").Value = 0.001001119613647461
$ws.Range("O39").Value = 0.0009989738464355469
$ws.Range("O40").Value = 0.06852245330810547
$ws.Range("O42").Value = 0.0744326114654541
$ws.Range("O43").Value = 0.002000808715820312
$ws.Range("O44").Value = 0.09596085548400879
$ws.Range("O46").Value = 0
$ws.Range("O47").Value = 0.007974863052368164
$ws.Range("O48").Value = 0.0776817798614502
$ws.Range("O49").Value = 0.001585721969604492
$ws.Range("O51").Value = 0.002999067306518555
$ws.Range("O52").Value = 0.002012014389038086
$ws.Range("O53").Value = 0.001000165939331055
$ws.Range("O54").Value = 0.002000093460083008
$ws.Range("O55").Value = 0.04021143913269043
$ws.Range("O56").Value = 0.002002716064453125
$ws.Range("O57").Value = 0.004666566848754883
$ws.Range("O58").Value = 0.002027034759521484
$ws.Range("O59").Value = 0
$ws.Range("O60").Value = 0.002007246017456055
$ws.Range("O61").Value = 0.001001358032226562
$ws.Range("O62").Value = 0.01867985725402832
$ws.Range("O64").Value = 0.001348972320556641
$ws.Range("O65").Value = 0.07170486450195312
$ws.Range("O66").Value = 0.01494503021240234
$ws.Range("O67").Value = 0
$ws.Range("O68").Value = 0.001005887985229492
$ws.Range("O69").Value = 0.008698225021362305
$ws.Range("O70").Value = 0.001046180725097656
$ws.Range("O71").Value = 0.02065706253051758
$ws.Range("O72").Value = 130.2574288845062
$ws.Range("O73").Value = 0.01621150970458984
$ws.Range("O75").Value = 0.001000404357910156
$ws.Range("O77").Value = 0.002464771270751953
$ws.Range("O78").Value = 0
$ws.Range("O79").Value = 0
$ws.Range("O80").Value = 0.001611471176147461
$ws.Range("O81").Value = 0.05254340171813965
$ws.Range("O82").Value = 0.002088308334350586
$ws.Range("O83").Value = 0
$ws.Range("O84").Value = 0.006002664566040039
$ws.Range("O85").Value = 0.01372933387756348
$ws.Range("O86").Value = 0.002716779708862305
$ws.Range("O87").Value = 0.004767417907714844
$ws.Range("O88").Value = 0.001913785934448242
$ws.Range("O89").Value = 0.002114534378051758
$ws.Range("O90").Value = 0.1540229320526123
$ws.Range("O92").Value = 0.003784656524658203
$ws.Range("O94").Value = 0.001026153564453125
$ws.Range("O95").Value = 0.001916170120239258
$ws.Range("O96").Value = 0.08431196212768555
$ws.Range("O97").Value = 0.002912282943725586
$ws.Range("O98").Value = 0.001999616622924805
$ws.Range("O99").Value = 0.001999616622924805
$ws.Range("O100").Value = 0.07430744171142578
$ws.Range("O102").Value = 0.001999616622924805
$ws.Range("O103").Value = 0.03700923919677734
$ws.Range("O104").Value = 4.572839736938477
$ws.Range("O105").Value = 0.03751778602600098
$ws.Range("O106").Value = 0.07468533515930176
$ws.Range("O107").Value = 0.05855035781860352
$ws.Range("O108").Value = 0.007000446319580078
$ws.Range("O109").Value = 0.002020597457885742
$ws.Range("O110").Value = 0.06022143363952637
$ws.Range("O111").Value = 0.002999067306518555
$ws.Range("O112").Value = 0.02746891975402832
$ws.Range("O113").Value = 0.03130412101745605
$ws.Range("O114").Value = 0.01950812339782715
$ws.Range("O115").Value = 0.01691293716430664
$ws.Range("O116").Value = 0.02046847343444824
$ws.Range("O117").Value = 0.01300621032714844
$ws.Range("O118").Value = 0.00571441650390625
$ws.Range("O119").Value = 0.001009941101074219
$ws.Range("O120").Value = 0.002621173858642578
$ws.Range("O121").Value = 0.00109100341796875
$ws.Range("O122").Value = 0.002629995346069336
$ws.Range("O123").Value = 0
$ws.Range("O124").Value = 0.001965045928955078
$ws.Range("O125").Value = 0.001003742218017578
$ws.Range("O126").Value = 0.0009987354278564453
$ws.Range("O127").Value = 0.001024246215820312
$ws.Range("O128").Value = 0.001600980758666992
$ws.Range("O129").Value = 0.001001834869384766
$ws.Range("O130").Value = 0.0009977817535400391
$ws.Range("O131").Value = 0.001990795135498047
$ws.Range("O132").Value = 0
$ws.Range("O135").Value = 0.001412391662597656
$ws.Range("O136").Value = 0.00101161003112793
$ws.Range("O139").Value = 0.002086400985717773
$ws.Range("O140").Value = 0.003625392913818359
$ws.Range("O141").Value = 0.002510547637939453
$ws.Range("O142").Value = 0
$ws.Range("O143").Value = 0.006009578704833984
$ws.Range("O146").Value = 0
$ws.Range("O148").Value = 0
$ws.Range("O150").Value = 0.2442600727081299
$ws.Range("O151").Value = 0.002309083938598633
$ws.Range("O153").Value = 0
$ws.Range("O155").Value = 0.01589012145996094
$ws.Range("O161").Value = 0
$ws.Range("O165").Value = 0.01001429557800293
$ws.Range("O167").Value = 0.003986597061157227
$ws.Range("O168").Value = 0
$ws.Range("O169").Value = 0.001318454742431641
$ws.Range("O170").Value = 0
$ws.Range("O172").Value = 0
$ws.Range("O175").Value = 0.005607843399047852
